$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1802.8889
$ws.Range("I80").Value = 659.6
$ws.Range("K80").Value = 1978.8
$ws.Range("M80").Value = -980.8000000000002
$ws.Range("H83").Value = 1802.8889
$ws.Range("I83").Value = 659.6
$ws.Range("K83").Value = 5936.400000000001
$ws.Range("M83").Value = -944.4000000000005
$ws.Range("H88").Value = 2723
$ws.Range("J88").Value = 2833.75
$ws.Range("L88").Value = 2833.75
$ws.Range("N88").Value = -3645.75
$ws.Range("H91").Value = 2723
$ws.Range("J91").Value = 2833.75
$ws.Range("L91").Value = 2833.75
$ws.Range("N91").Value = -5641.75
$ws.Range("H98").Value = 5153.923
$ws.Range("I98").Value = 4634.8696
$ws.Range("K98").Value = 4634.8696
$ws.Range("M98").Value = -3136.8696
$ws.Range("H122").Value = 5153.923
$ws.Range("I122").Value = 4634.8696
$ws.Range("K122").Value = 13904.6088
$ws.Range("M122").Value = -11454.6088
$ws.Range("H132").Value = 4582.846
$ws.Range("J132").Value = 4250
$ws.Range("L132").Value = 12750
$ws.Range("N132").Value = -17810

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 8464.25
$ws.Range("I88").Value = 8528.444
$ws.Range("J88").Value = 8433.842000000001
$ws.Range("K88").Value = 8528.444
$ws.Range("L88").Value = 8433.842000000001
$ws.Range("M88").Value = -8122.444
$ws.Range("N88").Value = -9245.842000000001
$ws.Range("H91").Value = 8464.25
$ws.Range("I91").Value = 8528.444
$ws.Range("J91").Value = 8433.842000000001
$ws.Range("K91").Value = 8528.444
$ws.Range("L91").Value = 8433.842000000001
$ws.Range("M91").Value = -7124.444
$ws.Range("N91").Value = -11241.842
$ws.Range("H122").Value = 1870.5
$ws.Range("I122").Value = 1208.8572
$ws.Range("K122").Value = 3626.5716
$ws.Range("M122").Value = -1176.5716
$ws.Range("H132").Value = 1915.9183
$ws.Range("I132").Value = 1797.409
$ws.Range("J132").Value = 2958.8
$ws.Range("K132").Value = 5392.227000000001
$ws.Range("L132").Value = 8876.400000000001
$ws.Range("M132").Value = -2862.227000000001
$ws.Range("N132").Value = -13936.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1043.2727
$ws.Range("J80").Value = 942.1111
$ws.Range("L80").Value = 942.1111
$ws.Range("N80").Value = -2938.1111
$ws.Range("H83").Value = 1043.2727
$ws.Range("J83").Value = 942.1111
$ws.Range("L83").Value = 4710.555499999999
$ws.Range("N83").Value = -14694.5555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4671.5835
$ws.Range("I31").Value = 2687
$ws.Range("K31").Value = 2687
$ws.Range("M31").Value = -2392
$ws.Range("H34").Value = 4671.5835
$ws.Range("I34").Value = 2687
$ws.Range("K34").Value = 2687
$ws.Range("M34").Value = -2485
$ws.Range("H58").Value = 3840.24
$ws.Range("J58").Value = 4609.0835
$ws.Range("L58").Value = 4609.0835
$ws.Range("N58").Value = -5015.0835
$ws.Range("H59").Value = 100999.4
$ws.Range("I59").Value = 14999
$ws.Range("J59").Value = 122499.5
$ws.Range("K59").Value = 14999
$ws.Range("L59").Value = 122499.5
$ws.Range("M59").Value = -13854
$ws.Range("N59").Value = -124789.5
$ws.Range("H99").Value = 3115.6843
$ws.Range("I99").Value = 2906.25
$ws.Range("J99").Value = 3268
$ws.Range("K99").Value = 2906.25
$ws.Range("L99").Value = 3268
$ws.Range("M99").Value = -1408.25
$ws.Range("N99").Value = -6264
$ws.Range("H126").Value = 3115.6843
$ws.Range("I126").Value = 2906.25
$ws.Range("J126").Value = 3268
$ws.Range("K126").Value = 8718.75
$ws.Range("L126").Value = 9804
$ws.Range("M126").Value = -6248.75
$ws.Range("N126").Value = -14744
$ws.Range("H130").Value = 45997.5
$ws.Range("J130").Value = 45997.5
$ws.Range("L130").Value = 45997.5
$ws.Range("N130").Value = -56037.5
$ws.Range("H132").Value = 2299.8235
$ws.Range("I132").Value = 2299.8235
$ws.Range("K132").Value = 6899.470499999999
$ws.Range("M132").Value = -4369.470499999999
$ws.Range("H134").Value = 6083.6875
$ws.Range("I134").Value = 5626.48
$ws.Range("J134").Value = 7716.5713
$ws.Range("K134").Value = 16879.44
$ws.Range("L134").Value = 23149.7139
$ws.Range("M134").Value = -14344.44
$ws.Range("N134").Value = -28219.7139
$ws.Range("H136").Value = 3840.24
$ws.Range("J136").Value = 4609.0835
$ws.Range("L136").Value = 13827.2505
$ws.Range("N136").Value = -18927.2505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2200.7058
$ws.Range("I132").Value = 1056.4445
$ws.Range("J132").Value = 3488
$ws.Range("K132").Value = 9508.0005
$ws.Range("L132").Value = 31392
$ws.Range("M132").Value = -6978.0005
$ws.Range("N132").Value = -36452

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 20009
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 20009
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 20009
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -21067
$ws.Range("H80").Value = 3417.1538
$ws.Range("J80").Value = 3442.3
$ws.Range("L80").Value = 3442.3
$ws.Range("N80").Value = -5438.3
$ws.Range("H83").Value = 3417.1538
$ws.Range("J83").Value = 3442.3
$ws.Range("L83").Value = 17211.5
$ws.Range("N83").Value = -27195.5
$ws.Range("H102").Value = 5491.9565
$ws.Range("I102").Value = 5238.8667
$ws.Range("K102").Value = 5238.8667
$ws.Range("M102").Value = -3616.8667
$ws.Range("H113").Value = 2890.25
$ws.Range("J113").Value = 2993.5
$ws.Range("L113").Value = 2993.5
$ws.Range("N113").Value = -7333.5
$ws.Range("H126").Value = 6857.143
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -5030
$ws.Range("H135").Value = 92150
$ws.Range("J135").Value = 92150
$ws.Range("L135").Value = 92150
$ws.Range("N135").Value = -102290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9792.6
$ws.Range("I16").Value = 11504.625
$ws.Range("K16").Value = 11504.625
$ws.Range("M16").Value = -11334.625
$ws.Range("H68").Value = 2280.1482
$ws.Range("J68").Value = 2412.818
$ws.Range("L68").Value = 2412.818
$ws.Range("N68").Value = -3910.818
$ws.Range("H71").Value = 2280.1482
$ws.Range("J71").Value = 2412.818
$ws.Range("L71").Value = 12064.09
$ws.Range("N71").Value = -19552.09
$ws.Range("H82").Value = 1980.2727
$ws.Range("I82").Value = 1749.8462
$ws.Range("J82").Value = 2313.111
$ws.Range("K82").Value = 1749.8462
$ws.Range("L82").Value = 2313.111
$ws.Range("M82").Value = -1388.8462
$ws.Range("N82").Value = -3035.111
$ws.Range("H85").Value = 1980.2727
$ws.Range("I85").Value = 1749.8462
$ws.Range("J85").Value = 2313.111
$ws.Range("K85").Value = 1749.8462
$ws.Range("L85").Value = 2313.111
$ws.Range("M85").Value = -501.8462
$ws.Range("N85").Value = -4809.111
$ws.Range("H122").Value = 5673.3105
$ws.Range("I122").Value = 4900
$ws.Range("J122").Value = 5700.9287
$ws.Range("K122").Value = 14700
$ws.Range("L122").Value = 17102.7861
$ws.Range("M122").Value = -12250
$ws.Range("N122").Value = -22002.7861
$ws.Range("H136").Value = 3749.875
$ws.Range("I136").Value = 2600
$ws.Range("K136").Value = 7800
$ws.Range("M136").Value = -5250

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 349.33334
$ws.Range("I113").Value = 349.33334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1048.00002
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1121.99998
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 316959.12
$ws.Range("I122").Value = 503970.1
$ws.Range("K122").Value = 1511910.3
$ws.Range("M122").Value = -1509460.3
$ws.Range("H132").Value = 2581.82
$ws.Range("I132").Value = 2314.9534
$ws.Range("J132").Value = 4221.143
$ws.Range("K132").Value = 6944.860199999999
$ws.Range("L132").Value = 12663.429
$ws.Range("M132").Value = -4414.860199999999
$ws.Range("N132").Value = -17723.429
